$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (GitHub Actions scheduled update).
# A new coin (WrappedBTC) now appears at row 17; every coin that used
# to follow it shifts down by one row, and the coin that used to be
# last (EnergySwap) drops off the bottom of the fixed 50-row table.
# Column A (the 0-based rank index) is untouched by this refresh.
#
# Price text looks like plain numbers (e.g. "209.24"), but the sheet
# stores it as literal text (note values like "26.332.76" or the
# subscript-notation "0.0₃0736" which are not valid numbers at all).
# Assigning such a numeric-looking string via .Value normally makes
# Excel coerce it to a real number, so for those cells we briefly force
# a Text format, assign the literal string, then restore the default
# "Normal" style so no visible formatting change is left behind.

$ws.Range("D2").Value = '26.332.76'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '1.585.80'
$ws.Range("E3").Value = '  -1.00%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0611'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.244'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.21%  '

$ws.Range("D12").Value = '1.809.40'
$ws.Range("E12").Value = '  -0.94%  '

$ws.Range("D13").Value = '1.588.96'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.515'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.57%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.321.98'
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  -1.86%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '205.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.24%  '

$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.17%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.81%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.19%  '

$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.94%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.49%  '

$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.10%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.29%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.09%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.00%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.28%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.40%  '

$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.279.58'
$ws.Range("E34").Value = '  -2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.17%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("E37").Value = '  -1.82%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.595'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.02%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0167'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.49%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.37%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.766'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.14%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.90%  '

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.721.69'
$ws.Range("E45").Value = '  -0.90%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.23%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.38%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.101'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₇0985'
$ws.Range("E50").Value = '  -6.65%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '
